$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESOURCES")

# Insert a new blank row above the current row 2 ("Natural Gas"), shifting the
# existing data rows (Natural Gas / Electricity / Solar) down by one.
$ws.Rows.Item(2).Insert()

# Seed the new row's formatting: start from the header row's look (this is
# what gives column F - the "reference" column - its border/center style),
# then overwrite columns A:E with the formatting used by the other data rows.
$ws.Range("A1:F1").Copy() | Out-Null
$ws.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A2:E2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new "none" resource row: no PEN/CO2/costs impact, no reference.
$ws.Range("A2").Value = "none"
$ws.Range("B2").Value = "NONE"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = ""

$ws.Range("A2:F2").Select()
